# Add a new "remilia_scarlet" (touhou) artist row to Sheet1, right before the
# existing "rennkuu" row (row 266), pushing everything else down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 266 (shifts rows 266.. down by one).
$ws.Rows.Item(266).Insert()

# Fill in the new row's data (style is inherited from the surrounding rows,
# same centered alignment already used throughout the table).
$ws.Cells.Item(266, 1).Value = "remilia_scarlet"
$ws.Cells.Item(266, 2).Value = "danbooru"
$ws.Cells.Item(266, 3).Value = "character"
$ws.Cells.Item(266, 4).Value = "touhou"

# Update the view state to match what was recorded after the edit: scrolled
# down so row 253 is the top visible row, with D267 as the active cell.
$excel.ActiveWindow.ScrollRow = 253
$ws.Range("D267").Select()
